$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the "preterm" column (AB) to a single value "Term" for all
# data rows (2-18), replacing the mix of "Preterm"/"Term" entries. This
# removes the now-unused "Preterm" shared string.
for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 28).Value = "Term"
}

# Fix up session_error / session_error_type notes for a few rows.
$ws.Range("H8").Value = "All orders were randomized"
$ws.Range("H9").Value = "All orders were randomized"
$ws.Range("G17").Value = "error"

# Column width adjustment: column K (11) gets its own width of 12,
# separate from the rest of columns L:S (12-19) which stay at 8.5.
# (ColumnWidth uses Excel's character-width units, which is 5/6 less
# than the stored OOXML width, so 12 - 5/6 = 11.1666... gives width=12.)
$ws.Range("K1").EntireColumn.ColumnWidth = 11.166666666666666

# Update the selected cell/range shown when the sheet is active.
$ws.Range("H21").Select()
